$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.257.77'
$ws.Range("E2").Value = '  +0.21%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.786.74'
$ws.Range("E3").Value = '  -0.27%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.011'
$ws.Range("E4").Value = '  +0.61%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '328.52'
$ws.Range("E5").Value = '  -2.77%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.005'
$ws.Range("E6").Value = '  +0.39%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4377'
$ws.Range("E7").Value = '  -3.28%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3740'
$ws.Range("E8").Value = '  +4.59%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.57'
$ws.Range("E9").Value = '  +0.05%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07592'
$ws.Range("E10").Value = '  +1.74%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.136'
$ws.Range("E11").Value = '  -0.15%  '

$ws.Range("B12").Value = 'BinanceUSD'
$ws.Range("C12").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.011'
$ws.Range("E12").Value = '  +0.74%  '

$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.40'
$ws.Range("E13").Value = '  +0.26%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.244'
$ws.Range("E14").Value = '  +0.62%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.409'
$ws.Range("E15").Value = '  +2.49%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.786.09'
$ws.Range("E16").Value = '  -0.38%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001087'
$ws.Range("E17").Value = '  +0.38%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06711'
$ws.Range("E18").Value = '  +0.33%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '81.92'
$ws.Range("E19").Value = '  +1.10%  '

$ws.Range("E20").Value = '  +0.30%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.42'
$ws.Range("E21").Value = '  +1.47%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.242'
$ws.Range("E22").Value = '  -2.15%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.354.83'
$ws.Range("E23").Value = '  +0.55%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.61'
$ws.Range("E24").Value = '  -1.77%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.439'
$ws.Range("E25").Value = '  +2.12%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '20.57'
$ws.Range("E26").Value = '  +0.96%  '

$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.366'
$ws.Range("E27").Value = '  -0.22%  '

$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '151.94'
$ws.Range("E28").Value = '  -1.08%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.988.22'
$ws.Range("E29").Value = '  -0.42%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.303'
$ws.Range("E30").Value = '  +3.04%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '130.84'
$ws.Range("E31").Value = '  -1.15%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.993'
$ws.Range("E32").Value = '  -1.95%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.843'
$ws.Range("E33").Value = '  -0.33%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.09260'
$ws.Range("E34").Value = '  -1.51%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.2240'
$ws.Range("E35").Value = '  +3.93%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.10'
$ws.Range("E36").Value = '  +0.30%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6672'
$ws.Range("E37").Value = '  +0.65%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06299'
$ws.Range("E38").Value = '  +1.30%  '

$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.225'
$ws.Range("E39").Value = '  +1.21%  '

$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02322'
$ws.Range("E40").Value = '  -1.78%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.207'
$ws.Range("E41").Value = '  -0.17%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.438'
$ws.Range("E42").Value = '  -2.94%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.059'
$ws.Range("E43").Value = '  +0.18%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.005'
$ws.Range("E44").Value = '  +0.49%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.95'
$ws.Range("E45").Value = '  +0.12%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6107'
$ws.Range("E46").Value = '  +0.88%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.812'
$ws.Range("E47").Value = '  -1.24%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '127.42'
$ws.Range("E48").Value = '  -0.54%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.018'
$ws.Range("E49").Value = '  +0.01%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06997'
$ws.Range("E50").Value = '  -1.21%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.139'
$ws.Range("E51").Value = '  -1.81%  '
